$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E1:E4").ClearContents()
